# Updates the cryptos list figures (price + 1h volume change) on the active
# worksheet to match the latest scrape, and swaps the TrustWalletToken /
# VeChain rows (42 and 43) which changed rank order.
#
# Note: several "Price" values are plain decimal-looking strings (e.g.
# "244.29") that must remain TEXT (as in the source data) rather than be
# auto-converted to numbers by Excel. We force those via the classic
# leading-apostrophe "text prefix" so Excel stores them as strings, exactly
# like the multi-dot values (e.g. "36.680.21") which Excel leaves as text
# on its own because they aren't valid numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.680.21'
$ws.Range("E2").Value = '  +0.39%  '

$ws.Range("D3").Value = '1.965.63'
$ws.Range("E3").Value = '  +1.15%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = "'244.29"
$ws.Range("E5").Value = '  +0.10%  '

$ws.Range("D6").Value = "'0.617"
$ws.Range("E6").Value = '  +0.52%  '

$ws.Range("D7").Value = "'58.46"
$ws.Range("E7").Value = '  +1.31%  '

$ws.Range("D9").Value = "'0.373"
$ws.Range("E9").Value = '  +1.28%  '

$ws.Range("E10").Value = '  -3.93%  '

$ws.Range("E11").Value = '  -0.11%  '

$ws.Range("D12").Value = "'22.20"
$ws.Range("E12").Value = '  +3.64%  '

$ws.Range("D13").Value = '2.254.95'
$ws.Range("E13").Value = '  +1.31%  '

$ws.Range("D14").Value = "'0.823"
$ws.Range("E14").Value = '  -0.35%  '

$ws.Range("E15").Value = '  +0.81%  '

$ws.Range("E16").Value = '  +0.30%  '

$ws.Range("D17").Value = '1.968.44'
$ws.Range("E17").Value = '  +0.25%  '

$ws.Range("D18").Value = '36.579.11'
$ws.Range("E18").Value = '  +0.39%  '

$ws.Range("D19").Value = "'69.71"
$ws.Range("E19").Value = '  +0.11%  '

$ws.Range("D20").Value = '0.0₃0860'
$ws.Range("E20").Value = '  -0.93%  '

$ws.Range("E21").Value = '  +1.65%  '

$ws.Range("D22").Value = "'228.54"
$ws.Range("E22").Value = '  -0.39%  '

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = '  -0.18%  '

$ws.Range("E24").Value = '  -3.14%  '

$ws.Range("E25").Value = '  +1.83%  '

$ws.Range("E26").Value = '  +0.67%  '

$ws.Range("D27").Value = "'0.138"
$ws.Range("E27").Value = '  +9.33%  '

$ws.Range("D28").Value = "'160.09"
$ws.Range("E28").Value = '  -1.55%  '

$ws.Range("D29").Value = "'19.39"
$ws.Range("E29").Value = '  -0.25%  '

$ws.Range("E30").Value = '  +1.09%  '

$ws.Range("D31").Value = "'1.12"
$ws.Range("E31").Value = '  -2.59%  '

$ws.Range("E32").Value = '  +0.49%  '

$ws.Range("E33").Value = '  -1.89%  '

$ws.Range("E34").Value = '  -1.01%  '

$ws.Range("E35").Value = '  +0.12%  '

$ws.Range("D36").Value = "'6.10"
$ws.Range("E36").Value = '  +0.52%  '

$ws.Range("D37").Value = "'3.41"
$ws.Range("E37").Value = '  +15.62%  '

$ws.Range("E38").Value = '  +4.06%  '

$ws.Range("E39").Value = '  -0.68%  '

$ws.Range("E40").Value = '  +2.81%  '

$ws.Range("E41").Value = '  -0.15%  '

# Row 42 and 43 swap: TrustWalletToken moves up, VeChain moves down.
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = "'1.17"
$ws.Range("E42").Value = '  -0.84%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = "'0.0212"
$ws.Range("E43").Value = '  +1.44%  '

$ws.Range("D44").Value = "'15.99"
$ws.Range("E44").Value = '  +0.15%  '

$ws.Range("D45").Value = '1.367.05'
$ws.Range("E45").Value = '  +1.17%  '

$ws.Range("E46").Value = '  +0.16%  '

$ws.Range("D47").Value = "'87.46"
$ws.Range("E47").Value = '  -0.35%  '

$ws.Range("E48").Value = '  -1.10%  '

$ws.Range("E49").Value = '  +0.73%  '

$ws.Range("D50").Value = '2.146.13'
$ws.Range("E50").Value = '  +1.32%  '

$ws.Range("D51").Value = "'43.42"
$ws.Range("E51").Value = '  -5.02%  '
